# Generate Report for Handback
# - Status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   on the Overview sheet (E/F) and both language sheets (col C).
# - Refresh "Latest Handback DateTime" (col K) on zh-cn / de-de sheets.
# - Clear the stale "version mismatch" error detail (col P) on zh-cn / de-de,
#   since the handback is now in sync.
# - Widen column C (status) on the language sheets / E,F on Overview to fit
#   the new, longer status text; narrow column P (error detail) since it's
#   now unused.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("K2").Value = "2016-10-18 13:21:20"
$zhcn.Range("K3").Value = "2016-10-18 13:21:20"

$zhcn.Range("P2").Value = ""

$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668
$zhcn.Columns.Item(16).ColumnWidth = 12.833333333333334

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("K2").Value = "2016-10-18 13:21:38"
$dede.Range("K3").Value = "2016-10-18 13:21:38"

$dede.Range("P2").Value = ""

$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(16).ColumnWidth = 12.833333333333334
